{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\n// Identify the three paragraphs that must be removed:\n//  1) the blank paragraph right after \"Janeiro: Editora Interci\u00eancia , 2004.\"\n//  2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//  3) \"\u00a9 2020 . Contact: luizeleno@usp.br. ... Creative Commons Attribution\"\n// They always appear as this exact three-paragraph run, so locate the\n// \"Ver no Jupiter\" paragraph and remove it together with its immediate\n// neighbors (the blank line before it and the copyright line after it).\nconst items = paragraphs.items;\nlet jupiterIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIndex = i;\n    break;\n  }\n}\n\nif (jupiterIndex !== -1) {\n  const toDelete = [];\n  if (jupiterIndex - 1 >= 0 && items[jupiterIndex - 1].text.trim() === \"\") {\n    toDelete.push(items[jupiterIndex - 1]);\n  }\n  toDelete.push(items[jupiterIndex]);\n  if (\n    jupiterIndex + 1 < items.length &&\n    items[jupiterIndex + 1].text.indexOf(\"Creative Commons Attribution\") !== -1\n  ) {\n    toDelete.push(items[jupiterIndex + 1]);\n  }\n\n  for (const para of toDelete) {\n    para.delete();\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph containing \"Ver no Jupiter\" (the footer line that\n# was removed from the page). The three paragraphs that must disappear\n# are: the blank line right before it, the \"Ver no Jupiter...\" line\n# itself, and the \"\u00a9 2020 ...\" copyright line right after it.\n$jupIdx = -1\n$i = 1\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Ver no Jupiter*\") {\n        $jupIdx = $i\n    }\n    $i = $i + 1\n}\n\nif ($jupIdx -ge 1) {\n    $startIdx = $jupIdx\n    if ($jupIdx - 1 -ge 1) {\n        $prevText = $d.Paragraphs.Item($jupIdx - 1).Range.Text\n        if ($prevText.Trim() -eq \"\") {\n            $startIdx = $jupIdx - 1\n        }\n    }\n\n    $endIdx = $jupIdx\n    if ($jupIdx + 1 -le $d.Paragraphs.Count) {\n        $nextText = $d.Paragraphs.Item($jupIdx + 1).Range.Text\n        if ($nextText -like \"*Creative Commons Attribution*\") {\n            $endIdx = $jupIdx + 1\n        }\n    }\n\n    $startPara = $d.Paragraphs.Item($startIdx)\n    $endPara = $d.Paragraphs.Item($endIdx)\n\n    $rangeToDelete = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $rangeToDelete.Delete()\n}\n"}
